$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Create the brand-new shared strings first, and in the exact order they
# were introduced, so sharedStrings.xml ends up with the same ordering.
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "教务处信息点击"
$ws.Range("A21").Value = "主页重新排版"
$ws.Range("C17").Value = [string][char]215
$ws.Range("E17").Value = "找不到合适的处理逻辑过程"
$ws.Range("A22").Value = "最小化按钮"
$ws.Range("A19").Value = "应用图标背景透明"

# ---------------------------------------------------------------------------
# Row 13: add C13 ("√"), D13 (date 42792), E13 (blank, centered style)
#         and change B13 from style s=3 to s=2 (center+vcenter date style)
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy($ws.Range("B13"))
$ws.Range("B13").Value = 42792
$ws.Range("C2").Copy($ws.Range("C13"))
$ws.Range("D2").Copy($ws.Range("D13"))
$ws.Range("D13").Value = 42792
$ws.Range("E13").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 14: change B14 style s=3 -> s=2 ; add blank centered C14, D14, E14
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy($ws.Range("B14"))
$ws.Range("B14").Value = 42792
$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("E14").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 15: change B15 style s=3 -> s=2 ; add blank centered C15, D15, E15
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy($ws.Range("B15"))
$ws.Range("B15").Value = 42792
$ws.Range("C15").HorizontalAlignment = -4108
$ws.Range("D15").HorizontalAlignment = -4108
$ws.Range("E15").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 16: change row height, B16/D16 style s=3 -> s=2 ; add blank E16
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).RowHeight = 33
$ws.Range("B2").Copy($ws.Range("B16"))
$ws.Range("B16").Value = 42792
$ws.Range("D2").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 42792
$ws.Range("E16").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 17: change row height, add B17 (date), style C17/leave text, D17 (blank),
#         style E17 text already set above.
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).RowHeight = 42.75
$ws.Range("B2").Copy($ws.Range("B17"))
$ws.Range("B17").Value = 42791
$ws.Range("A2").Copy($ws.Range("C17"))
$ws.Range("C17").Value = [string][char]215
$ws.Range("D17").HorizontalAlignment = -4108
$ws.Range("A2").Copy($ws.Range("E17"))
$ws.Range("E17").Value = "找不到合适的处理逻辑过程"

# ---------------------------------------------------------------------------
# Row 18: set row height (was none), B18/D18 style s=3 -> s=2 ; add blank E18
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).RowHeight = 31.5
$ws.Range("B2").Copy($ws.Range("B18"))
$ws.Range("B18").Value = 42792
$ws.Range("D2").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 42792
$ws.Range("E18").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 19 (new): 应用图标背景透明 (text already set above)
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 34.5
$ws.Range("A2").Copy($ws.Range("A19"))
$ws.Range("A19").Value = "应用图标背景透明"
$ws.Range("B2").Copy($ws.Range("B19"))
$ws.Range("B19").Value = 42792
$ws.Range("C2").Copy($ws.Range("C19"))
$ws.Range("D2").Copy($ws.Range("D19"))
$ws.Range("D19").Value = 42793
$ws.Range("E19").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 20 (new): 教务处信息点击 (text already set above)
# ---------------------------------------------------------------------------
$ws.Rows.Item(20).RowHeight = 33.75
$ws.Range("A2").Copy($ws.Range("A20"))
$ws.Range("A20").Value = "教务处信息点击"
$ws.Range("B2").Copy($ws.Range("B20"))
$ws.Range("B20").Value = 42792
$ws.Range("C2").Copy($ws.Range("C20"))
$ws.Range("D20").HorizontalAlignment = -4108
$ws.Range("E20").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 21 (new): 主页重新排版 (text already set above)
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).RowHeight = 35.25
$ws.Range("A2").Copy($ws.Range("A21"))
$ws.Range("A21").Value = "主页重新排版"
$ws.Range("B2").Copy($ws.Range("B21"))
$ws.Range("B21").Value = 42792
$ws.Range("C21").HorizontalAlignment = -4108
$ws.Range("D21").HorizontalAlignment = -4108
$ws.Range("E21").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 22 (new): 最小化按钮 (text already set above)
# ---------------------------------------------------------------------------
$ws.Rows.Item(22).RowHeight = 42
$ws.Range("A2").Copy($ws.Range("A22"))
$ws.Range("A22").Value = "最小化按钮"

# ---------------------------------------------------------------------------
# Column E width (engine quantizes to 1/7 character units; 25.5 is exactly
# between two reachable values, use the nearer/round-to-even one)
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 24.75

# ---------------------------------------------------------------------------
# Selection / view
# ---------------------------------------------------------------------------
[void]$ws.Range("E18").Select()

Write-Output "done"
